$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 502.4762
$ws.Range("I32").Value = 436.5
$ws.Range("J32").Value = 543.0769
$ws.Range("K32").Value = 436.5
$ws.Range("L32").Value = 543.0769
$ws.Range("M32").Value = -110.5
$ws.Range("N32").Value = -1195.0769

$ws.Range("H112").Value = 10490660
$ws.Range("I112").Value = 795.5
$ws.Range("J112").Value = 11364816
$ws.Range("K112").Value = 2386.5
$ws.Range("L112").Value = 34094448
$ws.Range("M112").Value = -1278.5
$ws.Range("N112").Value = -34096664

$ws.Range("H129").Value = 1059.871
$ws.Range("J129").Value = 1085.3
$ws.Range("L129").Value = 3255.9
$ws.Range("N129").Value = -13255.9

$ws.Range("H138").Value = 6573024
$ws.Range("I138").Value = 7942284.5
$ws.Range("J138").Value = 6415032.5
$ws.Range("K138").Value = 23826853.5
$ws.Range("L138").Value = 19245097.5
$ws.Range("M138").Value = -23821713.5
$ws.Range("N138").Value = -19255377.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 3840
$ws.Range("I88").Value = 10000
$ws.Range("J88").Value = 2300
$ws.Range("K88").Value = 10000
$ws.Range("L88").Value = 2300
$ws.Range("M88").Value = -9594
$ws.Range("N88").Value = -3112

$ws.Range("H91").Value = 3840
$ws.Range("I91").Value = 10000
$ws.Range("J91").Value = 2300
$ws.Range("K91").Value = 10000
$ws.Range("L91").Value = 2300
$ws.Range("M91").Value = -8596
$ws.Range("N91").Value = -5108

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1756.6428
$ws.Range("I86").Value = 2057.8
$ws.Range("J86").Value = 1589.3334
$ws.Range("K86").Value = 2057.8
$ws.Range("L86").Value = 1589.3334
$ws.Range("M86").Value = -934.8000000000002
$ws.Range("N86").Value = -3835.3334

$ws.Range("H89").Value = 1756.6428
$ws.Range("I89").Value = 2057.8
$ws.Range("J89").Value = 1589.3334
$ws.Range("K89").Value = 10289
$ws.Range("L89").Value = 7946.666999999999
$ws.Range("M89").Value = -4673
$ws.Range("N89").Value = -19178.667

$ws.Range("H134").Value = 4500827.5
$ws.Range("I134").Value = 846538.6
$ws.Range("J134").Value = 23816354
$ws.Range("K134").Value = 2539615.8
$ws.Range("L134").Value = 71449062
$ws.Range("M134").Value = -2537080.8
$ws.Range("N134").Value = -71454132

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 125450
$ws.Range("I4").Value = 514.2857
$ws.Range("J4").Value = 1000000
$ws.Range("K4").Value = 514.2857
$ws.Range("L4").Value = 1000000
$ws.Range("M4").Value = -402.2857
$ws.Range("N4").Value = -1000224

$ws.Range("H7").Value = 124.166664
$ws.Range("I7").Value = 124.166664
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 124.166664
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -11.166664
$ws.Range("N7").ClearContents()

$ws.Range("H31").Value = 6547.593
$ws.Range("I31").Value = 3933
$ws.Range("J31").Value = 7854.8887
$ws.Range("K31").Value = 3933
$ws.Range("L31").Value = 7854.8887
$ws.Range("M31").Value = -3638
$ws.Range("N31").Value = -8444.8887

$ws.Range("H34").Value = 6547.593
$ws.Range("I34").Value = 3933
$ws.Range("J34").Value = 7854.8887
$ws.Range("K34").Value = 3933
$ws.Range("L34").Value = 7854.8887
$ws.Range("M34").Value = -3731
$ws.Range("N34").Value = -8258.8887

$ws.Range("H99").Value = 166668180
$ws.Range("I99").Value = 3000
$ws.Range("K99").Value = 3000
$ws.Range("M99").Value = -1502

$ws.Range("H126").Value = 166668180
$ws.Range("I126").Value = 3000
$ws.Range("K126").Value = 9000
$ws.Range("M126").Value = -6530

$ws.Range("H141").Value = 302607.4
$ws.Range("J141").Value = 310481.75
$ws.Range("L141").Value = 310481.75
$ws.Range("N141").Value = -320841.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 457.88235
$ws.Range("I114").Value = 123
$ws.Range("J114").Value = 692.3
$ws.Range("K114").Value = 369
$ws.Range("L114").Value = 2076.9
$ws.Range("M114").Value = 2885
$ws.Range("N114").Value = -8584.9

$ws.Range("H117").Value = 444.14285
$ws.Range("I117").Value = 221.8
$ws.Range("K117").Value = 665.4000000000001
$ws.Range("M117").Value = 2776.6

$ws.Range("H129").Value = 1319.5217
$ws.Range("I129").Value = 425.375
$ws.Range("J129").Value = 1796.4
$ws.Range("K129").Value = 1276.125
$ws.Range("L129").Value = 5389.200000000001
$ws.Range("M129").Value = 3723.875
$ws.Range("N129").Value = -15389.2

$ws.Range("H131").Value = 4567655
$ws.Range("I131").Value = 750
$ws.Range("J131").Value = 4763379.5
$ws.Range("K131").Value = 2250
$ws.Range("L131").Value = 14290138.5
$ws.Range("M131").Value = 2790
$ws.Range("N131").Value = -14300218.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1716.2273
$ws.Range("I97").Value = 1641.625
$ws.Range("J97").Value = 1915.1666
$ws.Range("K97").Value = 1641.625
$ws.Range("L97").Value = 1915.1666
$ws.Range("M97").Value = -1145.625
$ws.Range("N97").Value = -2907.1666

$ws.Range("H122").Value = 4269.16
$ws.Range("I122").Value = 3642.5293
$ws.Range("J122").Value = 5600.75
$ws.Range("K122").Value = 10927.5879
$ws.Range("L122").Value = 16802.25
$ws.Range("M122").Value = -8477.5879
$ws.Range("N122").Value = -21702.25

$ws.Range("H136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3599.85
$ws.Range("I7").Value = 3499.25
$ws.Range("J7").Value = 3625
$ws.Range("K7").Value = 3499.25
$ws.Range("L7").Value = 3625
$ws.Range("M7").Value = -3387.25
$ws.Range("N7").Value = -3849

$ws.Range("H82").Value = 0
$ws.Range("I82").Value = 0
$ws.Range("K82").Value = 0
$ws.Range("M82").ClearContents()

$ws.Range("H85").Value = 0
$ws.Range("I85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("M85").ClearContents()

$ws.Range("H122").Value = 3383.6843
$ws.Range("I122").Value = 2381.6667
$ws.Range("J122").Value = 3846.1538
$ws.Range("K122").Value = 7145.000100000001
$ws.Range("L122").Value = 11538.4614
$ws.Range("M122").Value = -4695.000100000001
$ws.Range("N122").Value = -16438.4614

$ws.Range("H126").Value = 3599.85
$ws.Range("I126").Value = 3499.25
$ws.Range("J126").Value = 3625
$ws.Range("K126").Value = 10497.75
$ws.Range("L126").Value = 10875
$ws.Range("M126").Value = -8027.75
$ws.Range("N126").Value = -15815

$ws.Range("H127").Value = 48500
$ws.Range("J127").Value = 48500
$ws.Range("L127").Value = 48500
$ws.Range("N127").Value = -58420

$ws.Range("H139").Value = 50715
$ws.Range("J139").Value = 50715
$ws.Range("L139").Value = 50715
$ws.Range("N139").Value = -60995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1470
$ws.Range("I122").Value = 1250
$ws.Range("J122").Value = 1800
$ws.Range("K122").Value = 3750
$ws.Range("L122").Value = 5400
$ws.Range("M122").Value = -1300
$ws.Range("N122").Value = -10300

$ws.Range("H123").Value = 35000
$ws.Range("J123").Value = 35000
$ws.Range("L123").Value = 35000
$ws.Range("N123").Value = -44800

$ws.Range("H132").Value = 2950.804
$ws.Range("I132").Value = 2889.7896
$ws.Range("J132").Value = 3129.1538
$ws.Range("K132").Value = 8669.3688
$ws.Range("L132").Value = 9387.4614
$ws.Range("M132").Value = -6139.3688
$ws.Range("N132").Value = -14447.4614
